# Update NATMI TPM-derived expression / specificity values in Fgf22-Fgfrl1 sheet
# (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> ECs
$ws.Range("I2").Value = 0.3046368955123587
$ws.Range("J2").Value = 0.3046368955123587
$ws.Range("M2").Value = 1.170631
$ws.Range("N2").Value = 2.341262
$ws.Range("O2").Value = 0.1596624636318675
$ws.Range("P2").Value = 0.1281983861842902
$ws.Range("Q2").Value = 0.03679566380233333
$ws.Range("R2").Value = 0.220773982814
$ws.Range("S2").Value = 0.048639077250667
$ws.Range("T2").Value = 0.03905395837687663

# Row 3: FAPs -> FAPs
$ws.Range("I3").Value = 0.3046368955123587
$ws.Range("J3").Value = 0.3046368955123587
$ws.Range("O3").Value = 0.4101137613801331
$ws.Range("P3").Value = 0.4939412918191532
$ws.Range("S3").Value = 0.12493578307374
$ws.Range("T3").Value = 0.1504727417051509

# Row 4: FAPs -> Inflammatory-Mac
$ws.Range("I4").Value = 0.3046368955123587
$ws.Range("J4").Value = 0.3046368955123587
$ws.Range("M4").Value = 0.1597873333333333
$ws.Range("N4").Value = 0.479362
$ws.Range("O4").Value = 0.02179340825346879
$ws.Range("P4").Value = 0.02624799565280337
$ws.Range("Q4").Value = 0.005022488723777778
$ws.Range("R4").Value = 0.04520239851400001
$ws.Range("S4").Value = 0.006639076232970148
$ws.Range("T4").Value = 0.007996107909091908

# Row 5: FAPs -> MuSCs
$ws.Range("I5").Value = 0.3046368955123587
$ws.Range("J5").Value = 0.3046368955123587
$ws.Range("M5").Value = 2.5622985
$ws.Range("N5").Value = 5.124597
$ws.Range("O5").Value = 0.349472114671693
$ws.Range("P5").Value = 0.2806029676494365
$ws.Range("Q5").Value = 0.08053902055149999
$ws.Range("R5").Value = 0.483234123309
$ws.Range("S5").Value = 0.1064621000817236
$ws.Range("T5").Value = 0.08548201693627917

# Row 6: FAPs -> Neutrophils
$ws.Range("I6").Value = 0.3046368955123587
$ws.Range("J6").Value = 0.3046368955123587
$ws.Range("M6").Value = 0.2687716666666666
$ws.Range("N6").Value = 0.8063149999999999
$ws.Range("O6").Value = 0.03665779093022745
$ws.Range("P6").Value = 0.04415066821064279
$ws.Range("Q6").Value = 0.008448120617222221
$ws.Range("R6").Value = 0.07603308555499999
$ws.Range("S6").Value = 0.01116731562532559
$ws.Range("T6").Value = 0.0134499224984864

# Row 7: FAPs -> Resolving-Mac
$ws.Range("I7").Value = 0.3046368955123587
$ws.Range("J7").Value = 0.3046368955123587
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.163505
$ws.Range("N7").Value = 0.490515
$ws.Range("O7").Value = 0.02230046113261011
$ws.Range("P7").Value = 0.02685869048367381
$ws.Range("Q7").Value = 0.005139343661666667
$ws.Range("R7").Value = 0.04625409295500001
$ws.Range("S7").Value = 0.006793543247932361
$ws.Range("T7").Value = 0.008182148086473724

# Row 8: Inflammatory-Mac -> ECs
$ws.Range("G8").Value = 0.07174733333333333
$ws.Range("H8").Value = 0.215242
$ws.Range("I8").Value = 0.6953631044876413
$ws.Range("J8").Value = 0.6953631044876413
$ws.Range("M8").Value = 1.170631
$ws.Range("N8").Value = 2.341262
$ws.Range("O8").Value = 0.1596624636318675
$ws.Range("P8").Value = 0.1281983861842902
$ws.Range("Q8").Value = 0.08398965256733333
$ws.Range("R8").Value = 0.503937915404
$ws.Range("S8").Value = 0.1110233863812005
$ws.Range("T8").Value = 0.08914442780741359

# Row 9: Inflammatory-Mac -> FAPs
$ws.Range("G9").Value = 0.07174733333333333
$ws.Range("H9").Value = 0.215242
$ws.Range("I9").Value = 0.6953631044876413
$ws.Range("J9").Value = 0.6953631044876413
$ws.Range("O9").Value = 0.4101137613801331
$ws.Range("P9").Value = 0.4939412918191532
$ws.Range("Q9").Value = 0.2157383241362222
$ws.Range("R9").Value = 1.941644917226
$ws.Range("S9").Value = 0.2851779783063931
$ws.Range("T9").Value = 0.3434685501140023

# Row 10: Inflammatory-Mac -> Inflammatory-Mac
$ws.Range("G10").Value = 0.07174733333333333
$ws.Range("H10").Value = 0.215242
$ws.Range("I10").Value = 0.6953631044876413
$ws.Range("J10").Value = 0.6953631044876413
$ws.Range("M10").Value = 0.1597873333333333
$ws.Range("N10").Value = 0.479362
$ws.Range("O10").Value = 0.02179340825346879
$ws.Range("P10").Value = 0.02624799565280337
$ws.Range("Q10").Value = 0.01146431506711111
$ws.Range("R10").Value = 0.103178835604
$ws.Range("S10").Value = 0.01515433202049865
$ws.Range("T10").Value = 0.01825188774371147

# Row 11: Inflammatory-Mac -> MuSCs
$ws.Range("G11").Value = 0.07174733333333333
$ws.Range("H11").Value = 0.215242
$ws.Range("I11").Value = 0.6953631044876413
$ws.Range("J11").Value = 0.6953631044876413
$ws.Range("M11").Value = 2.5622985
$ws.Range("N11").Value = 5.124597
$ws.Range("O11").Value = 0.349472114671693
$ws.Range("P11").Value = 0.2806029676494365
$ws.Range("Q11").Value = 0.183838084579
$ws.Range("R11").Value = 1.103028507474
$ws.Range("S11").Value = 0.2430100145899694
$ws.Range("T11").Value = 0.1951209507131574

# Row 12: Inflammatory-Mac -> Neutrophils
$ws.Range("G12").Value = 0.07174733333333333
$ws.Range("H12").Value = 0.215242
$ws.Range("I12").Value = 0.6953631044876413
$ws.Range("J12").Value = 0.6953631044876413
$ws.Range("M12").Value = 0.2687716666666666
$ws.Range("N12").Value = 0.8063149999999999
$ws.Range("O12").Value = 0.03665779093022745
$ws.Range("P12").Value = 0.04415066821064279
$ws.Range("Q12").Value = 0.01928365035888889
$ws.Range("R12").Value = 0.17355285323
$ws.Range("S12").Value = 0.02549047530490186
$ws.Range("T12").Value = 0.03070074571215638

# Row 13: Inflammatory-Mac -> Resolving-Mac
$ws.Range("G13").Value = 0.07174733333333333
$ws.Range("H13").Value = 0.215242
$ws.Range("I13").Value = 0.6953631044876413
$ws.Range("J13").Value = 0.6953631044876413
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.163505
$ws.Range("N13").Value = 0.490515
$ws.Range("O13").Value = 0.02230046113261011
$ws.Range("P13").Value = 0.02685869048367381
$ws.Range("Q13").Value = 0.01173104773666667
$ws.Range("R13").Value = 0.10557942963
$ws.Range("S13").Value = 0.01550691788467774
$ws.Range("T13").Value = 0.01867654239720009
